$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "23-03-2025"
$ws.Range("B3").Value = "Sunrisers Hyderabad vs Rajasthan Royals"
$ws.Range("C3").Value = "Rajasthan Royals"
$ws.Range("D3").Value = "Rajasthan Royals"

$ws.Range("A4").Value = "23-03-2025"
$ws.Range("B4").Value = "Chennai Super Kings vs Mumbai Indians"
$ws.Range("C4").Value = "Mumbai Indians"
$ws.Range("D4").Value = "Mumbai Indians"
